$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6829547882080078
$ws.Range("B1").Value = 2.016331195831299
$ws.Range("C1").Value = 4.806789398193359
$ws.Range("D1").Value = 1.837946772575378
$ws.Range("E1").Value = 1.219586253166199
